$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet uses "X" marks to flag which ODS category applies to each row.
# Several marks were moved to a different column (re-classification of the
# "ODS" table), and the newly-placed D4 mark also got a (no-visual-change)
# font format applied to it.

# Row 2: mark moved from C2 -> D2
$ws.Range("D2").Value = "X"
$ws.Range("C2").ClearContents()

# Row 3: mark moved from B3 -> C3
$ws.Range("C3").Value = "X"
$ws.Range("B3").ClearContents()

# Row 4: mark moved from C4 -> D4, with a font format applied
$ws.Range("D4").Value = "X"
$ws.Range("C4").ClearContents()
$ws.Range("D4").Font.Bold = $true

# Row 6: mark moved from D6 -> B6
$ws.Range("B6").Value = "X"
$ws.Range("D6").ClearContents()

# Update the selected/active cell to the newly edited cell
$ws.Range("D4").Select() | Out-Null

# Configure printing: paper size 9 = A4, orientation 1 = portrait
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
